# Apply the workbook edit described by the commit:
# "Clean up directories - move 2018.2_AR71791 to 2018.2, add lime/shared, fix gaussian_delay.py."
#
# The functional spreadsheet change is on Sheet1 ("Emulation Control" table):
#   - C28 (PL0 / IOPLL, "Divisor 0")  changes from 6 to 8
#   - C30 (PL1 / IOPLL, "Divisor 0")  changes from 5 to 8
# All other cells in the sheet are formulas that depend (directly or
# transitively) on these two inputs, so Excel will recompute them
# automatically once the values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two divisor input cells that drive the rest of the sheet.
$ws.Range("C28").Value = 8
$ws.Range("C30").Value = 8

# Recalculate all open workbooks so dependent formulas refresh immediately.
$excel.CalculateFullRebuild()

# Reflect the cell selection that was active when the workbook was saved.
$ws.Activate()
$ws.Range("C30").Select()
